$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Rafa, 2025-09-29, 0.3
$ws.Range("A9").Value = "Rafa"
$ws.Range("B9").Value = 45929
$ws.Range("C9").Value = 0.3

# Row 10: mike, 2025-09-30, 0.4
$ws.Range("A10").Value = "mike"
$ws.Range("B10").Value = 45930
$ws.Range("C10").Value = 0.4

# Apply the same date number format already used by the column above (d-mmm)
$ws.Range("B9:B10").NumberFormatLocal = "d-mmm"

$ws.Range("A11").Select()
